$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K2 (row 2): keep the AVERAGEIF formula, switch to the "0.00"
#     number-format (same style already used by L2); value is
#     recalculated automatically because row 9 below gains new inputs.
$ws.Range("K2").NumberFormat = "0.00"

# --- Row 3: F3 / G3 become running-total formulas ----------------------
$ws.Range("F3").Clear()
$ws.Range("F3").Formula = "=C2+C3"
$ws.Range("G3").Clear()
$ws.Range("G3").Formula = "=D2+D3"

# --- Row 5: F5 / G5 become formulas; H5 date added ----------------------
$ws.Range("F5").Clear()
$ws.Range("F5").Formula = "=C4+C5"
$ws.Range("G5").Clear()
$ws.Range("G5").Formula = "=D4+D5"
$ws.Range("H5").Value = 45859
$ws.Range("H5").NumberFormat = "dd/mm/yyyy"

# --- Row 7: F7 / G7 become formulas -------------------------------------
$ws.Range("F7").Clear()
$ws.Range("F7").Formula = "=C6+C7"
$ws.Range("G7").Clear()
$ws.Range("G7").Formula = "=D6+D7"

# --- Row 9: brand-new data for 2020 H2 + running totals + date ---------
$ws.Range("C9").Value = 353
$ws.Range("D9").Value = 400
$ws.Range("E9").Value = 1
$ws.Range("F9").Formula = "=C8+C9"
$ws.Range("G9").Formula = "=D8+D9"
$ws.Range("H9").Value = 45875
$ws.Range("H9").NumberFormat = "dd/mm/yyyy"

# --- Row 11: F11 / G11 become formulas ----------------------------------
$ws.Range("F11").Clear()
$ws.Range("F11").Formula = "=C10+C11"
$ws.Range("G11").Clear()
$ws.Range("G11").Formula = "=D10+D11"

# --- Row 17: F17 / G17 become formulas ----------------------------------
$ws.Range("F17").Clear()
$ws.Range("F17").Formula = "=C16+C17"
$ws.Range("G17").Clear()
$ws.Range("G17").Formula = "=D16+D17"

# --- Row 19: F19 / G19 become formulas ----------------------------------
$ws.Range("F19").Clear()
$ws.Range("F19").Formula = "=C18+C19"
$ws.Range("G19").Clear()
$ws.Range("G19").Formula = "=D18+D19"

# --- Row heights: rows that now carry a running-total formula settle to
#     the tighter 13.8pt height already seen elsewhere on the sheet.
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8
$ws.Rows.Item(11).RowHeight = 13.8
$ws.Rows.Item(17).RowHeight = 13.8
$ws.Rows.Item(19).RowHeight = 13.8
